$xlPasteFormats = -4122
$xlPortrait = 1

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# Add the new worksheet at the end of the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "addAbsenceTest"

# Fill the sheet column by column (matches the shared-string insertion order).
$newSheet.Range("A1").Value = "emp"
$newSheet.Range("A2").Value = "Twentyone 21"

$newSheet.Range("B1").Value = "absenceType"
$newSheet.Range("B2").Value = "Fiscal Year"

$newSheet.Range("C1").Value = "absenceStDt"
$newSheet.Range("C2").Value = "20/08/2022"

$newSheet.Range("D1").Value = "absenceEndDt"
$newSheet.Range("D2").Value = "20/08/2022"

$newSheet.Range("E1").Value = "absenceStTime"
$newSheet.Range("E2").Value = "10:00 AM"

$newSheet.Range("F1").Value = "absenceEndTime"
$newSheet.Range("F2").Value = "05:30 PM"

$newSheet.Range("G1").Value = "absenceLocation"
$newSheet.Range("G2").Value = "Cramer Dentistry"

$newSheet.Range("H1").Value = "absenceJob"
$newSheet.Range("H2").Value = "Design"

$newSheet.Range("I1").Value = "absenceNotes"
$newSheet.Range("I2").Value = "Adding the absence through automated test case"

# Reuse the existing "label" style (Consolas font) from addTimeSheetTest!A2.
$ws2.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial($xlPasteFormats)
$newSheet.Range("B2").PasteSpecial($xlPasteFormats)
$newSheet.Range("G2").PasteSpecial($xlPasteFormats)
$newSheet.Range("H2").PasteSpecial($xlPasteFormats)

# Reuse the existing quote-prefix style from addTimeSheetTest!B2.
$ws2.Range("B2").Copy()
$newSheet.Range("C2").PasteSpecial($xlPasteFormats)
$newSheet.Range("D2").PasteSpecial($xlPasteFormats)
$newSheet.Range("E2").PasteSpecial($xlPasteFormats)
$newSheet.Range("F2").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# Match addTimeSheetTest's page setup (portrait orientation).
$newSheet.PageSetup.Orientation = $xlPortrait

# The new sheet becomes the active/visible tab; the previous active sheet
# (addTimeSheetTest) loses its tabSelected flag.
$newSheet.Range("E7").Select()
$newSheet.Activate()
